$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Story Tasks" column (D) previously showed the same generic hyperlink
# text "Story Tasks Document" for rows C1-C5 (rows 3-7). Each row now shows
# its own distinct link label "Cx - Story Document" matching its story ID.
$ws.Range("D3").Value = "C1 - Story Document"
$ws.Range("D4").Value = "C2 - Story Document"
$ws.Range("D5").Value = "C3 - Story Document"
$ws.Range("D6").Value = "C4 - Story Document"
$ws.Range("D7").Value = "C5 - Story Document"
